# Auto-update draw results: append the 2025-10-14 Pick 4 row (row 28).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A28:E28")

# Force the new cells to Text format before writing so values that look like
# numbers/dates (the date "2025-10-14" and the phase code "251014") are
# stored as literal text strings, matching the rest of the column, instead
# of being auto-coerced into a date serial / number by Excel.
$newRow.NumberFormat = "@"

$ws.Range("A28").Value = "2025-10-14"
$ws.Range("B28").Value = "Pick 4"
$ws.Range("C28").Value = "251014"
$ws.Range("D28").Value = "7-8-7-6"
$ws.Range("E28").Value = "2025-10-14T21:36:47.194+04:00"

# Drop the temporary number format so the new row carries the same (default)
# style as every other data row instead of leaving a "Text" style behind.
$newRow.ClearFormats()
